$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.980.20'
$ws.Range('E2').Value = '  -4.46%  '
$ws.Range('D3').Value = '2.713.94'
$ws.Range('E3').Value = '  -6.81%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '492.54'
$ws.Range('E5').Value = '  -7.01%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '133.09'
$ws.Range('E6').Value = '  -7.09%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.521'
$ws.Range('E8').Value = '  -5.80%  '
$ws.Range('D9').Value = '2.717.05'
$ws.Range('E9').Value = '  -6.59%  '
$ws.Range('E10').Value = '  -0.34%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.100'
$ws.Range('E11').Value = '  -7.04%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.338'
$ws.Range('E12').Value = '  -4.08%  '
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('D14').Value = '3.199.86'
$ws.Range('E14').Value = '  -6.39%  '
$ws.Range('D15').Value = '58.140.25'
$ws.Range('E15').Value = '  -4.35%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '20.88'
$ws.Range('E16').Value = '  -7.86%  '
$ws.Range('D17').Value = '2.718.72'
$ws.Range('E17').Value = '  -6.69%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.0000131'
$ws.Range('E18').Value = '  -6.63%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.60'
$ws.Range('E19').Value = '  -6.65%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.66'
$ws.Range('E20').Value = '  -7.65%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '338.23'
$ws.Range('E21').Value = '  -6.39%  '
$ws.Range('E22').Value = '  -6.57%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.61'
$ws.Range('E24').Value = '  -1.17%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '61.62'
$ws.Range('E25').Value = '  -2.80%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.416'
$ws.Range('E26').Value = '  -7.60%  '
$ws.Range('E27').Value = '  -7.55%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('E29').Value = '  -6.57%  '
$ws.Range('D30').Value = '0.0₃0782'
$ws.Range('E30').Value = '  -9.17%  '
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('E32').Value = '  -6.34%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '18.66'
$ws.Range('E33').Value = '  -5.24%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '146.49'
$ws.Range('E34').Value = '  -5.11%  '
$ws.Range('E35').Value = '  -7.19%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.20'
$ws.Range('E36').Value = '  -7.05%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.888'
$ws.Range('E37').Value = '  -11.94%  '
$ws.Range('E38').Value = '  -9.05%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '36.24'
$ws.Range('E39').Value = '  -4.49%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.30%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '2.135.00'
$ws.Range('E41').Value = '  -8.56%  '
$ws.Range('E42').Value = '  -6.81%  '
$ws.Range('E43').Value = '  -4.14%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.31'
$ws.Range('E44').Value = '  -10.82%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.587'
$ws.Range('E45').Value = '  -8.59%  '
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '18.41'
$ws.Range('E47').Value = '  -12.01%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0220'
$ws.Range('E48').Value = '  -5.49%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0874'
$ws.Range('E49').Value = '  -5.38%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '4.45'
$ws.Range('E50').Value = '  -7.97%  '
$ws.Range('E51').Value = '  -7.86%  '
